$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" = strikeouts) values regenerated from actual K counts
# instead of the previous "Strike#" placeholder data.
$kValues = @{
    2 = 2
    3 = 3
    4 = 1
    5 = 0
    6 = 1
    7 = 0
    8 = 2
    9 = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 2
    22 = 2
    23 = 2
    24 = 2
    25 = 0
    26 = 2
    27 = 3
    28 = 1
    29 = 3
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 0
    35 = 2
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    40 = 2
    41 = 1
    42 = 1
    43 = 0
    44 = 1
    45 = 3
    46 = 0
    47 = 0
    48 = 1
    49 = 2
    50 = 0
    51 = 1
    52 = 2
    53 = 3
    54 = 1
    55 = 2
    56 = 1
    57 = 2
    58 = 2
    59 = 3
    60 = 2
    61 = 1
    63 = 1
    64 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

